$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 05:08:09"
$wsZhCn.Range("H2").Value = "2016-03-20 05:08:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 05:08:17"
$wsDeDe.Range("H2").Value = "2016-03-20 05:09:01"
